# Append new log rows (12-21) to the "logs" sheet, matching the
# 001 -> 004 result update described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Static column values re-used on every new row.
$colB = "12 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, b#max-digit-skip-all-punctuation >= 7, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii, b#ascii >= 6"
$colC = "11 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, b#max-digit-skip-all-punctuation >= 7, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii"
$colD = "12 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, b#max-digit-skip-all-punctuation >= 7, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii, b#digit >= 7"
$colE = "Neural-Network"

$rule1000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$rule2000 = "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"

$rows = @(
    @{ Row=12; Time="20160426_134527"; Model=$rule1000; Classify=0.914191419141914; Segment=0.47 }
    @{ Row=13; Time="20160426_140017"; Model=$rule1000; Classify=0.914191419141914; Segment=0.46 }
    @{ Row=14; Time="20160426_141632"; Model=$rule1000; Classify=0.914191419141914; Segment=0.53 }
    @{ Row=15; Time="20160426_143259"; Model=$rule1000; Classify=0.914191419141914; Segment=0.58 }
    @{ Row=16; Time="20160426_144837"; Model=$rule1000; Classify=0.914191419141914; Segment=0.57 }
    @{ Row=17; Time="20160426_150633"; Model=$rule2000; Classify=0.914191419141914; Segment=0.46 }
    @{ Row=18; Time="20160426_153658"; Model=$rule2000; Classify=0.914191419141914; Segment=0.46 }
    @{ Row=19; Time="20160426_160736"; Model=$rule2000; Classify=0.907590759075908; Segment=0.45 }
    @{ Row=20; Time="20160426_163805"; Model=$rule2000; Classify=0.904290429042904; Segment=0.44 }
    @{ Row=21; Time="20160426_170854"; Model=$rule2000; Classify=0.914191419141914; Segment=0.46 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.Time
    $ws.Cells.Item($rowNum, 2).Value = $colB
    $ws.Cells.Item($rowNum, 3).Value = $colC
    $ws.Cells.Item($rowNum, 4).Value = $colD
    $ws.Cells.Item($rowNum, 5).Value = $colE
    $ws.Cells.Item($rowNum, 6).Value = $r.Model
    $ws.Cells.Item($rowNum, 7).Value = $colE
    $ws.Cells.Item($rowNum, 8).Value = $r.Model
    $ws.Cells.Item($rowNum, 9).Value = $colE
    $ws.Cells.Item($rowNum, 10).Value = $r.Model
    $ws.Cells.Item($rowNum, 11).Value = $r.Classify
    $ws.Cells.Item($rowNum, 12).Value = $r.Segment
}
